# improving user iteraction and main()
#
# Rebuilds the endeca_attributes table: the old "Sales Order Hold" style
# attribute list (15 data rows, rows 2-15) is replaced with the new
# Endeca "Sales Order Line" attribute list (12 data rows, rows 2-13).
# The display_name column (E) is emptied out for every row; a handful of
# the now-blank E cells keep distinctive font formatting that used to sit
# on other rows of the old table, so that formatting is copied over
# before the old text is wiped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$xlPasteFormats = -4122

# --- 1. Preserve the two special font styles that exist in the sheet
#        before we start rewriting row contents ---
# Style "Arial 11" currently lives on E6; copy it (format only) onto the
# rows that will need it in the new table: E5, E10, E12.
$ws.Range("E6").Copy()
$ws.Range("E5").PasteSpecial($xlPasteFormats)
$ws.Range("E10").PasteSpecial($xlPasteFormats)
$ws.Range("E12").PasteSpecial($xlPasteFormats)

# Style "Arial 10 / black" currently lives on E14; copy it onto E13.
$ws.Range("E14").Copy()
$ws.Range("E13").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# --- 2. Overwrite rows 2-13 with the new attribute rows ---
$data = @(
    @{ Row = 2;  Attr = "LINE_ID";               Type = "mdex:string";   Profile = 1 },
    @{ Row = 3;  Attr = "SHIPMENT_PRIORITY_CODE"; Type = "mdex:string";   Profile = 1 },
    @{ Row = 4;  Attr = "ORDERED_ITEM";           Type = "mdex:string";   Profile = 1 },
    @{ Row = 5;  Attr = "ORDERED_QUANTITY";       Type = "mdex:int";      Profile = 2 },
    @{ Row = 6;  Attr = "PROMISE_DATE";           Type = "mdex:dateTime"; Profile = 2 },
    @{ Row = 7;  Attr = "SCHEDULE_SHIP_DATE";     Type = "mdex:dateTime"; Profile = 2 },
    @{ Row = 8;  Attr = "REQUEST_DATE";           Type = "mdex:dateTime"; Profile = 2 },
    @{ Row = 9;  Attr = "SHIPPED_QUANTITY";       Type = "mdex:int";      Profile = 2 },
    @{ Row = 10; Attr = "ACTUAL_SHIPMENT_DATE";   Type = "mdex:dateTime"; Profile = 2 },
    @{ Row = 11; Attr = "FLOW_STATUS_CODE";       Type = "mdex:string";   Profile = 1 },
    @{ Row = 12; Attr = "SET_NAME";               Type = "mdex:string";   Profile = 1 },
    @{ Row = 13; Attr = "TYPE_CODE";              Type = "mdex:string";   Profile = 1 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = 14
    $ws.Cells.Item($r, 2).Value = $item.Attr
    $ws.Cells.Item($r, 3).Value = $item.Type
    $ws.Cells.Item($r, 4).Value = $item.Profile
    # display_name column is no longer populated for any row
    $ws.Cells.Item($r, 5).ClearContents()
}

# Rows 6 and 11 no longer carry any special formatting now that they hold
# no display_name text - wipe the inherited formatting completely.
$ws.Range("E6").Clear()
$ws.Range("E11").Clear()

# --- 3. Remove the two rows that are no longer part of the table ---
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(14).Delete()

# --- 4. Update the selected cell to match the new interaction state ---
$ws.Range("B4").Select()
